$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts). Regenerated values below replace the
# previous Strike# derived figures, per the accompanying data regeneration.
$kValues = @{
    2 = 2
    3 = 2
    4 = 1
    5 = 1
    6 = 2
    7 = 2
    8 = 2
    9 = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 3
    14 = 3
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 2
    24 = 2
    25 = 2
    26 = 1
    27 = 3
    28 = 1
    29 = 1
    30 = 3
    31 = 0
    32 = 0
    33 = 2
    34 = 0
    35 = 2
    36 = 1
    37 = 1
    38 = 0
    39 = 0
    40 = 1
    41 = 2
    42 = 2
    43 = 1
    44 = 1
    45 = 0
    46 = 2
    47 = 4
    48 = 1
    49 = 1
    50 = 1
    51 = 2
    52 = 3
    53 = 2
    54 = 2
    55 = 3
    56 = 1
    57 = 3
    58 = 2
    59 = 2
    60 = 4
    61 = 3
    62 = 3
    63 = 1
    64 = 1
    65 = 2
    66 = 0
    67 = 1
    68 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

